# Apply the updated crypto price / 1h-volume data (scheduled GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.160.15"
$ws.Range("E2").Value = "  +3.17%  "
$ws.Range("D3").Value = "2.996.24"
$ws.Range("E3").Value = "  +3.00%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'564.51"
$ws.Range("E5").Value = "  +2.91%  "
$ws.Range("D6").Value = "'138.55"
$ws.Range("E6").Value = "  +10.54%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").Value = "'0.519"
$ws.Range("E8").Value = "  +3.39%  "
$ws.Range("D9").Value = "2.989.05"
$ws.Range("E9").Value = "  +3.09%  "
$ws.Range("E10").Value = "  +8.29%  "
$ws.Range("D11").Value = "'5.09"
$ws.Range("E11").Value = "  +8.97%  "
$ws.Range("E12").Value = "  +5.02%  "
$ws.Range("E13").Value = "  +9.40%  "
$ws.Range("D14").Value = "'33.71"
$ws.Range("E14").Value = "  +4.56%  "
$ws.Range("E15").Value = "  +2.68%  "
$ws.Range("D16").Value = "3.491.92"
$ws.Range("E16").Value = "  +3.14%  "
$ws.Range("E17").Value = "  +8.03%  "
$ws.Range("D18").Value = "2.990.19"
$ws.Range("E18").Value = "  +3.09%  "
$ws.Range("D19").Value = "59.171.39"
$ws.Range("E19").Value = "  +3.29%  "
$ws.Range("D20").Value = "'428.48"
$ws.Range("E20").Value = "  +5.75%  "
$ws.Range("D21").Value = "'13.56"
$ws.Range("E21").Value = "  +6.00%  "
$ws.Range("D22").Value = "'0.716"
$ws.Range("E22").Value = "  +6.80%  "
$ws.Range("D23").Value = "'13.49"
$ws.Range("D24").Value = "'7.09"
$ws.Range("E24").Value = "  +3.92%  "
$ws.Range("D25").Value = "'80.62"
$ws.Range("E25").Value = "  +3.97%  "
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("E27").Value = "  +0.24%  "
$ws.Range("E28").Value = "  +11.14%  "
$ws.Range("D29").Value = "'2.54"
$ws.Range("E29").Value = "  +3.35%  "
$ws.Range("D30").Value = "'7.72"
$ws.Range("E30").Value = "  +6.96%  "
$ws.Range("E31").Value = "  +4.11%  "
$ws.Range("E32").Value = "  +2.88%  "
$ws.Range("D33").Value = "'0.0985"
$ws.Range("E33").Value = "  -0.12%  "
$ws.Range("D34").Value = "'0.997"
$ws.Range("E34").Value = "  +9.03%  "
$ws.Range("D35").Value = "0.0₃0771"
$ws.Range("E35").Value = "  +22.78%  "
$ws.Range("D36").Value = "'5.80"
$ws.Range("E36").Value = "  +6.99%  "
$ws.Range("E37").Value = "  +3.11%  "
$ws.Range("D38").Value = "'49.01"
$ws.Range("E38").Value = "  +1.74%  "
$ws.Range("D39").Value = "'8.66"
$ws.Range("E39").Value = "  +5.23%  "
$ws.Range("E40").Value = "  +12.32%  "
$ws.Range("D41").Value = "'397.86"
$ws.Range("E41").Value = "  +9.86%  "
$ws.Range("E42").Value = "  +4.01%  "
$ws.Range("D43").Value = "2.747.19"
$ws.Range("E43").Value = "  +4.66%  "
$ws.Range("D44").Value = "'0.107"
$ws.Range("E44").Value = "  +1.17%  "
$ws.Range("E45").Value = "  +10.18%  "
$ws.Range("E46").Value = "  -0.02%  "
$ws.Range("D47").Value = "'122.77"
$ws.Range("E47").Value = "  +2.56%  "
$ws.Range("E48").Value = "  +2.47%  "
$ws.Range("E49").Value = "  +3.87%  "
$ws.Range("D50").Value = "'23.42"
$ws.Range("E50").Value = "  +3.03%  "
$ws.Range("D51").Value = "'32.18"
$ws.Range("E51").Value = "  +18.28%  "
